# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns for zh-cn / de-de + generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 13:03:58"

# zh-cn sheet: status + latest handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 13:03:54"

# de-de sheet: status (shares the same "2016-09-03 13:03:58" handoff date
# string as the Overview sheet)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 13:03:58"

# Widen the status columns slightly to fit the new, longer text
# ("Ready for handoff" is wider than "In Translation").
$newStatusColWidth = 16.333333333333336
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
